$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = "2024-06-15 15:13:04"
$ws.Range("D36").Value = 200
$ws.Range("E36").Value = 15

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 2
$ws.Range("C37").Value = "2024-06-15 15:13:05"
$ws.Range("D37").Value = 200
$ws.Range("E37").Value = 0
